$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 4 and row 5 values for columns A, B, E, F, G, Q, R
$cols = @("A", "B", "E", "F", "G", "Q", "R")
foreach ($col in $cols) {
    $addr4 = "$col" + "4"
    $addr5 = "$col" + "5"
    $v4 = $ws.Range($addr4).Value2
    $v5 = $ws.Range($addr5).Value2
    $ws.Range($addr4).Value2 = $v5
    $ws.Range($addr5).Value2 = $v4
}

# Swap row 6 and row 7 values for columns A, Q, R
$cols2 = @("A", "Q", "R")
foreach ($col in $cols2) {
    $addr6 = "$col" + "6"
    $addr7 = "$col" + "7"
    $v6 = $ws.Range($addr6).Value2
    $v7 = $ws.Range($addr7).Value2
    $ws.Range($addr6).Value2 = $v7
    $ws.Range($addr7).Value2 = $v6
}
